# Adds a new daily record (day 24 of June/2025) to the faturamento_diario
# sheet. The existing rows for May, April and March (previously rows 25-116)
# are pushed down by one row, and the new record is placed at row 25,
# right after the last existing June row (row 24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the old row 25 so every row from 25 to 116 shifts
# down by one (rows 25-117 afterwards).
$ws.Rows.Item(25).Insert()

# Fill in the newly inserted row with the new data point.
$ws.Cells.Item(25, 1).Value = 24
$ws.Cells.Item(25, 2).Value = 22602.82
$ws.Cells.Item(25, 3).Value = 6
$ws.Cells.Item(25, 4).Value = 2025
$ws.Cells.Item(25, 5).Value = "06/2025"
